# Formed the consolidated report
# Recompute the "Absent" column (H) from the "Real" column (E) for each
# attendance row: Absent = 1 - Real.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 3) { $lastRow = 21 }

for ($r = 3; $r -le $lastRow; $r++) {
    $realCell = $ws.Cells.Item($r, 5)   # Column E - Real
    $absentCell = $ws.Cells.Item($r, 8) # Column H - Absent

    $realValue = $realCell.Value()
    if ($null -eq $realValue -or $realValue -eq "") {
        continue
    }

    $absentCell.Value = 1 - [double]$realValue
}
